# Correcting Relevance Markers Walker (2018) - Wolters (2018)
# Update row 3 (metrics_sim_with_priors.json) values in the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.7257217847769029
$ws.Range("D3").Value = 0.9081364829396326
$ws.Range("E3").Value = 0.9803149606299213
$ws.Range("H3").Value = 0.5804237726098191
$ws.Range("I3").Value = 0.07900748670671727
$ws.Range("J3").Value = 0.6259842519685039
$ws.Range("K3").Value = 4143.283464566929

$ws.Range("Q3").Value = 36
$ws.Range("R3").Value = 225
$ws.Range("S3").Value = 1371
$ws.Range("T3").Value = 4885
$ws.Range("U3").Value = 10425
$ws.Range("V3").Value = 47577
$ws.Range("W3").Value = 47388
$ws.Range("X3").Value = 46242
$ws.Range("Y3").Value = 42728
$ws.Range("Z3").Value = 37188

$ws.Range("AF3").Value = 0.999244
$ws.Range("AG3").Value = 0.995274
$ws.Range("AH3").Value = 0.971205
$ws.Range("AI3").Value = 0.897402
$ws.Range("AJ3").Value = 0.781047
